# cooprov_software.pptx edits
#  - Shape 1 ("Surface station" box): shrink height from 5006700 EMU to 4916400 EMU
#  - Shape 2 ("Raspberry pi on Robot" box): shrink height from 5006700 EMU to 4916400 EMU
#  - Shape 2 text: "Raspberry pi on Robot" -> "RaspberryPi on Robot"
#
# NB: PowerPoint's COM surface works in points (1 pt = 12700 EMU) for Shape
# geometry. 4916400 EMU = 387.118110... pt; 387.11815 is comfortably inside
# the (empirically verified) band of point values that this host converts
# back to exactly 4916400 EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$surfaceStation = $s.Shapes.Item(1)
$surfaceStation.Height = 387.11815

$raspberryPi = $s.Shapes.Item(2)
$raspberryPi.Height = 387.11815
$raspberryPi.TextFrame.TextRange.Text = "RaspberryPi on Robot"
